# Fruta / hortaliza, semanal
# Rows 2-11 get their Fecha (D), Volumen (M), Precio minimo (N),
# Precio maximo (O), Precio promedio ponderado (P) and Precio $/Kg (S)
# values rotated/permuted across rows, per source data update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current values for the columns that change, keyed by row number.
$cols = @("D", "M", "N", "O", "P", "S")
$orig = @{}
for ($r = 2; $r -le 11; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $orig[$r] = $rowVals
}

# Mapping: new row number -> source (old) row number that its new values come from.
$map = @{
    2  = 7
    3  = 5
    4  = 2
    5  = 3
    6  = 4
    7  = 8
    8  = 10
    9  = 11
    10 = 6
    11 = 9
}

foreach ($newRow in $map.Keys) {
    $oldRow = $map[$newRow]
    $src = $orig[$oldRow]
    foreach ($c in $cols) {
        $ws.Range("$c$newRow").Value2 = $src[$c]
    }
}
